# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.180.41'
$ws.Range('D3').Value = '1.601.87'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '303.13'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3781'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '52.02'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +3.57%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3619'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.267'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.57%  '
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.76'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.594'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.412'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001244'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.48%  '
$ws.Range('D17').Value = '1.599.42'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '93.96'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.94%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06886'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('E20').Value = '  -1.48%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.536'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.10%  '
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('E23').Value = '  -1.03%  '
$ws.Range('D24').Value = '23.180.20'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.395'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.46%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.990'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +6.51%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '21.23'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.25%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '149.90'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.248'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '133.80'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.362'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.752'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.90%  '
$ws.Range('D33').Value = '1.778.86'
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.9651'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07481'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.35%  '
$ws.Range('E36').Value = '  -2.69%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02720'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2520'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.49%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.08798'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.16%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.075'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.61%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.7096'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.360'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.73%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '12.49'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.32%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '15.62'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.62%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6537'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.58%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.311'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.83%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.018'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '132.04'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.07947'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.203'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.95%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.202'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.26%  '
